$wb = $excel.ActiveWorkbook

# --- Sheet "data1" (sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("data1")

# Remove rows 10-14 (subjects 6-10 of sample 1) - data now ends at row 9
$ws1.Rows("10:14").Delete() | Out-Null

# --- Sheet "data2" (sheet2.xml) ---
$ws2 = $wb.Worksheets.Item("data2")

# Renumber the subject column (C5:C9) from 11..15 down to 6..10
$ws2.Range("C5").Value = 6
$ws2.Range("C6").Value = 7
$ws2.Range("C7").Value = 8
$ws2.Range("C8").Value = 9
$ws2.Range("C9").Value = 10

# Remove rows 10-14 (old subjects 16-20) - data now ends at row 9
$ws2.Rows("10:14").Delete() | Out-Null

# Update selection to match the recorded end-state
$ws2.Range("A10:XFD14").Select() | Out-Null

# Re-activate "data1" last so it remains the selected/active tab,
# and set its selection to match the recorded end-state.
$ws1.Activate() | Out-Null
$ws1.Range("F23").Select() | Out-Null
